$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Alaska (row 3): status moves from "Needs spider" to "Blocked", and a new
# comment is added noting the site gave no contents. ---
$ws.Range("D3").Value2 = "Blocked"
$ws.Range("H3").Value2 = "Can't get any contents"

# --- Connecticut (row 8): work is finished. Status -> "Done", the daughter/son
# search notes are updated with the final results, and a conservative estimate
# formula (875/43) is added in column G, matching the style used by the other
# "Conservative Estimate" cells (copy the number format from G17). ---
$ws.Range("D8").Value2 = "Done"
$ws.Range("E8").Value2 = "Found 27 with DAUGHTER* & 19 with *DAUGHTER. Copied easily to CSV"
$ws.Range("F8").Value2 = "Found 679 with son* and 2954 with *son. Was able to relatively easy copy/paste in xls and then dedup in R"

$ws.Range("G17").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Formula = "=875/43"

# --- Restore the view to the top of the sheet (frozen header row, selection on
# B1) instead of the previous scroll position near the bottom. ---
$ws.Range("B1").Select()
